# Apply updated CircadiPy cosinor analysis values (sawtooth_05 results)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.003381731363802154
$ws.Range("H2").Value = 0.009328914107040426
$ws.Range("K2").Value = 4.588900240994006
$ws.Range("L2").Value = "[1.358375016508969, 7.819425465479043]"
$ws.Range("M2").Value = 0.005558440902627915
$ws.Range("N2").Value = 0.005558440902627915
$ws.Range("O2").Value = -1.371105502467618
$ws.Range("P2").Value = "[-2.251631971942234, -0.49057903299300154]"
$ws.Range("Q2").Value = 0.002408516219151613
$ws.Range("R2").Value = 0.002408516219151613
$ws.Range("S2").Value = 11.44183781978351
$ws.Range("T2").Value = "[9.580492463744527, 13.303183175822497]"
$ws.Range("W2").Value = 5.671491491491629
$ws.Range("X2").Value = 2.029249249249302
$ws.Range("Y2").Value = 9.313733733733956

# Row 3
$ws.Range("E3").Value = 25.14000000000049
$ws.Range("G3").Value = [double]"3.72624563649282e-07"
$ws.Range("H3").Value = [double]"5.185536667207232e-06"
$ws.Range("K3").Value = 6.380136102588835
$ws.Range("L3").Value = "[3.204317569495023, 9.555954635682646]"
$ws.Range("M3").Value = [double]"9.690274129736665e-05"
$ws.Range("N3").Value = 0.0001938054825947333
$ws.Range("O3").Value = -2.402579366709311
$ws.Range("P3").Value = "[-2.8680005005744658, -1.9371582328441566]"
$ws.Range("S3").Value = 12.02899566431611
$ws.Range("T3").Value = "[10.416180591321213, 13.641810737311014]"
$ws.Range("W3").Value = 9.61309309309328
$ws.Range("X3").Value = 7.750870870871021
$ws.Range("Y3").Value = 11.47531531531554
